# FIX: Control library updates
# Update the first acceleration data point (Sheet1!B2) from 1 to 2,
# and leave the selection on B3 (matching the author's last-saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data edit: Sheet1!B2 1 -> 2
$ws.Range("B2").Value = 2

# Move/restore the active selection to B3
[void]$ws.Range("B3").Select()
